# Updated symbol list on Sat Jan 14 14:39:14 UTC 2023 with GitHub Actions
# Applies the latest Coinranking price/volume(1h) refresh to the "cryptos" sheet.
# Values are written as literal text (leading "'" forces text, matching the
# original inlineStr/text cell type instead of letting Excel auto-parse the
# numeric-looking strings and "%" suffixed strings into Number values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.21"
$ws.Range("E2").Value = "'6.99%"
$ws.Range("E3").Value = "'8.84%"
$ws.Range("D4").Value = "'5.279"
$ws.Range("E4").Value = "'3.93%"
$ws.Range("D5").Value = "'0.07520"
$ws.Range("E5").Value = "'11.84%"
$ws.Range("D6").Value = "'7.879"
$ws.Range("E6").Value = "'7.62%"
$ws.Range("D7").Value = "'3.755"
$ws.Range("E7").Value = "'9.30%"
$ws.Range("D8").Value = "'1.483"
$ws.Range("E8").Value = "'6.49%"
$ws.Range("D9").Value = "'0.9129"
$ws.Range("E9").Value = "'1.37%"
$ws.Range("D10").Value = "'0.01774"
$ws.Range("E10").Value = "'2,648.88%"
$ws.Range("D11").Value = "'0.1702"
$ws.Range("E11").Value = "'7.94%"
$ws.Range("D12").Value = "'0.07785"
$ws.Range("E12").Value = "'13.06%"
$ws.Range("D13").Value = "'0.08057"
$ws.Range("E13").Value = "'6.45%"
$ws.Range("D14").Value = "'0.03016"
$ws.Range("E14").Value = "'2.96%"
$ws.Range("E15").Value = "'9.98%"
$ws.Range("D16").Value = "'0.001488"
$ws.Range("E16").Value = "'-6.00%"
$ws.Range("D17").Value = "'0.04544"
$ws.Range("E17").Value = "'1.35%"
$ws.Range("D18").Value = "'0.006150"
$ws.Range("E18").Value = "'-6.31%"
$ws.Range("D19").Value = "'3.478"
$ws.Range("E19").Value = "'0.84%"
$ws.Range("D20").Value = "'2.232"
$ws.Range("E20").Value = "'0.07%"
$ws.Range("D21").Value = "'0.3309"
$ws.Range("E21").Value = "'3.22%"
$ws.Range("E22").Value = "'1.24%"
$ws.Range("D23").Value = "'4.454"
$ws.Range("E23").Value = "'10.04%"
$ws.Range("D24").Value = "'0.1635"
$ws.Range("E24").Value = "'3.53%"
$ws.Range("D25").Value = "'0.001211"
$ws.Range("E25").Value = "'0.72%"
$ws.Range("D26").Value = "'0.004438"
$ws.Range("D27").Value = "'0.0001387"
$ws.Range("E27").Value = "'18.73%"
$ws.Range("D28").Value = "'0.0001737"
$ws.Range("E28").Value = "'7.56%"
$ws.Range("D40").Value = "'0.04549"
$ws.Range("E40").Value = "'7.28%"
$ws.Range("D41").Value = "'0.007011"
$ws.Range("E41").Value = "'3.57%"
$ws.Range("D42").Value = "'0.1342"
$ws.Range("E42").Value = "'8.42%"
$ws.Range("D43").Value = "'0.002218"
$ws.Range("E43").Value = "'1.48%"
$ws.Range("D44").Value = "'0.01334"
$ws.Range("E44").Value = "'16.26%"
$ws.Range("D45").Value = "'0.00006221"
$ws.Range("E45").Value = "'8.55%"
$ws.Range("D46").Value = "'1.873"
$ws.Range("E46").Value = "'-3.34%"
$ws.Range("D47").Value = "'0.01297"
$ws.Range("E47").Value = "'-13.57%"
